$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 272, pushing existing rows (272..374) down to (273..375)
$ws.Rows(272).Insert()

# Populate the newly inserted row with the new weekly record
$ws.Range("A272").Value = 10
$ws.Range("B272").Value = "Vega Modelo de Temuco"
$ws.Range("C272").Value = "La Araucanía"
$ws.Range("D272").Value = 45229
$ws.Range("E272").Value = 9
$ws.Range("F272").Value = 100114007
$ws.Range("G272").Value = "Jengibre"
$ws.Range("H272").Value = "Sin especificar"
$ws.Range("I272").Value = "Primera"
$ws.Range("J272").Value = 25
$ws.Range("K272").Value = 35000
$ws.Range("L272").Value = 35000
$ws.Range("M272").Value = 35000
$ws.Range("N272").Value = "$/caja 13 kilos"
$ws.Range("O272").Value = "Perú"
$ws.Range("P272").Value = 2692
$ws.Range("Q272").Value = 13
$ws.Range("R272").Value = "Hortaliza"
